$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (approx, engine snaps to 1/6-character pixel grid)
$ws.Columns.Item(1).ColumnWidth = 14.68
$ws.Columns.Item(2).ColumnWidth = 13.84

# Update cell values A1:B32
$ws.Range("A1").Value = -0.020878416066977934
$ws.Range("B1").Value = 0.020759560837007029
$ws.Range("A2").Value = 0.025350785661276021
$ws.Range("B2").Value = -0.02584891362287145
$ws.Range("A3").Value = 0.14420539239945995
$ws.Range("B3").Value = -0.14477250556808485
$ws.Range("A4").Value = -0.13921899882362609
$ws.Range("B4").Value = 0.13868824679678227
$ws.Range("A5").Value = -0.13268824719741623
$ws.Range("B5").Value = 0.13164444565043532
$ws.Range("A6").Value = -0.072189297427079691
$ws.Range("B6").Value = 0.072119563724296842
$ws.Range("A7").Value = -0.052119564208579661
$ws.Range("B7").Value = 0.051986631932479099
$ws.Range("A8").Value = -0.031986632419918948
$ws.Range("B8").Value = 0.031917473744955238
$ws.Range("A9").Value = -0.025917474163428267
$ws.Range("B9").Value = 0.025871811481354534
$ws.Range("A10").Value = -0.019871811902142156
$ws.Range("B10").Value = 0.019871277413592736
$ws.Range("A11").Value = -0.051464427334636298
$ws.Range("B11").Value = 0.051387439587543327
$ws.Range("A12").Value = -0.045387440010140612
$ws.Range("B12").Value = 0.045148356811655965
$ws.Range("A13").Value = -0.039148357240953224
$ws.Range("B13").Value = 0.039083272816099246
$ws.Range("A14").Value = -0.027083273278099007
$ws.Range("B14").Value = 0.027051817796007072
$ws.Range("A15").Value = -0.021051818228520425
$ws.Range("B15").Value = 0.021027053565493503
$ws.Range("A16").Value = -0.015027053999269402
$ws.Range("B16").Value = 0.015004380588413913
$ws.Range("A17").Value = -0.0090043810238844557
$ws.Range("B17").Value = 0.0089999995487390905
$ws.Range("A18").Value = -0.036110346919908665
$ws.Range("B18").Value = 0.03609664634712928
$ws.Range("A19").Value = -0.027096646763957288
$ws.Range("B19").Value = 0.027013786157357611
$ws.Range("A20").Value = -0.018013786577666835
$ws.Range("B20").Value = 0.018004306053232
$ws.Range("A21").Value = -0.0090043064740124024
$ws.Range("B21").Value = 0.0089999995788776488
$ws.Range("A22").Value = -0.10935647919237113
$ws.Range("B22").Value = 0.10892790324463064
$ws.Range("A23").Value = -0.084623246490974857
$ws.Range("B23").Value = 0.084124584206148434
$ws.Range("A24").Value = -0.042124584802618159
$ws.Range("B24").Value = 0.04199999940049004
$ws.Range("A25").Value = -0.053455148634526495
$ws.Range("B25").Value = 0.053413877931603793
$ws.Range("A26").Value = -0.047413878343451898
$ws.Range("B26").Value = 0.047366809847783742
$ws.Range("A27").Value = -0.026681355154075881
$ws.Range("B27").Value = 0.026554361228508405
$ws.Range("A28").Value = -0.02055436164380442
$ws.Range("B28").Value = 0.02047969270159733
$ws.Range("A29").Value = -0.0084796931491322169
$ws.Range("B29").Value = 0.0084580457861243019
$ws.Range("A30").Value = 0.011541953724970355
$ws.Range("B30").Value = -0.011628689229765676
$ws.Range("A31").Value = 0.026628688767504016
$ws.Range("B31").Value = -0.026677827822735978
$ws.Range("A32").Value = -0.0060008285705634989
$ws.Range("B32").Value = 0.0059999995845165799
